# "Add files via upload" - re-upload of the inventory sheet with a handful
# of data corrections made in Excel before saving:
#   - Several "Desktop" type devices (rows 17,18,20,21) were mislabeled as
#     "Notebook" in TipoDispositivo (column H) - fixed to "Desktop".
#   - The model name for the M93p ThinkCentre desktops was entered
#     inconsistently in a few rows ("Desktop M93p (ThinkCentre) - Tipo 10AA")
#     vs. the standard text used everywhere else
#     ("M93p Desktop (ThinkCentre) - Type 10AA"); normalized to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix TipoDispositivo (column H) for the V520s / Optiplex desktop rows ---
$ws.Range("H17").Value = "Desktop"
$ws.Range("H18").Value = "Desktop"
$ws.Range("H20").Value = "Desktop"
$ws.Range("H21").Value = "Desktop"

# --- Normalize ModeloDispositivo (column B) text for the M93p rows ---
$modelo = "M93p Desktop (ThinkCentre) - Type 10AA"
$ws.Range("B28").Value = $modelo
$ws.Range("B30").Value = $modelo
$ws.Range("B37").Value = $modelo
$ws.Range("B38").Value = $modelo

# --- Re-fit column B now that its contents were edited ---
$ws.Columns("B").AutoFit()

# --- Restore the view/selection state saved with the workbook ---
$ws.Range("H20:H21").Select()
